$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that would otherwise be auto-parsed as numbers
# (mirrors the source inlineStr cell type in the workbook XML)
$textCells = @("D5", "D10", "D11", "D15", "D16", "D19", "D21", "D22", "D23", "D24", "D25", "D29", "D30", "D32", "D33", "D36", "D38", "D39", "D43", "D45", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data
$ws.Range("D2").Value = "26.778.31"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.648.87"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "215.02"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "19.32"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.878.44"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "1.639.08"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "0.532"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "65.99"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "26.809.06"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "218.95"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "6.36"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").Value = "9.47"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  +9.30%  "
$ws.Range("D25").Value = "147.80"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").Value = "0.0520"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Value = "1.275.19"
$ws.Range("E34").Value = "  +9.10%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("E37").Value = "  +4.20%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.809"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.515"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "5.36"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "1.787.77"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "93.79"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("D47").Value = "56.00"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "7.67"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  +3.16%  "
